$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B3: change from text "3" to a true number 3
$ws.Range("B3").Value = 3

# Add new row 4 data
$ws.Range("A4").Value = "Ruilin"

# B4 must stay a text value "2" (not a number). Entering it with a leading
# apostrophe forces Excel to store it as text; then re-apply a plain
# (unstyled) style so no stray "quote prefix" number format is left behind.
$ws.Range("B4").Value = "'2"
$ws.Range("B4").Style = $ws.Range("A4").Style

$ws.Range("C4").Value = "does not provide any insight beyond, a strong reject"
$ws.Range("D4").Value = "SMY"
$ws.Range("E4").Value = "OTH"
$ws.Range("F4").Value = "1a2deef4-16ae-43c8-afd3-8fd2e076505e"
$ws.Range("G4").Value = "rJr4kfWCb_annotated.xlsx"
$ws.Range("H4").Value = "Overall, the paper does not provide any insight beyond: i tried this, i tried that and this works better than that; a strong reject."
